# mail-showcase.xlsx edit:
#  - web system-variable list (#system!V) gains a new entry "dragTo(fromLocator,xOffset,yOffset)"
#    inserted just above the existing list (shifts V63:V119 down to V64:V120).
#  - xml system-variable list (#system!AA) gains two new entries "beautify(xml,var)" and
#    "minify(xml,var)" inserted just above the existing list (shifts AA9:AA11 down to AA11:AA13).
#  - the two named ranges "web" and "xml" grow to cover the extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- shift column V (the "web" list) down by one row, starting at the bottom ---
for ($r = 119; $r -ge 63; $r--) {
    $val = $ws.Cells.Item($r, 22).Value2
    $ws.Cells.Item($r + 1, 22).Value = $val
}
$ws.Cells.Item(63, 22).Value = "dragTo(fromLocator,xOffset,yOffset)"

# --- shift column AA (the "xml" list) down by two rows, starting at the bottom ---
for ($r = 11; $r -ge 9; $r--) {
    $val = $ws.Cells.Item($r, 27).Value2
    $ws.Cells.Item($r + 2, 27).Value = $val
}
$ws.Cells.Item(9, 27).Value = "beautify(xml,var)"
$ws.Cells.Item(10, 27).Value = "minify(xml,var)"

# --- grow the named ranges to match the new list extents ---
$wb.Names.Item("web").RefersTo = "='#system'!`$V`$2:`$V`$120"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AA`$2:`$AA`$13"
